$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell contents to the new v0.4 structure:
# B2 now carries the option payload as JSON (was bare "json")
$ws.Range("B2").Value = '{"type":"json"}'
# B3 now carries the file path wrapped in a JSON object (was the bare path).
# The path's backslashes must be doubled so the cell text is valid JSON.
$ws.Range("B3").Value = '{"file":"C:\\Users\\xihu_\\Desktop\\2.json"}'

# Move the active selection from D1 to E5
[void]$ws.Range("E5").Select()
